# "Update Notes and Writeoff"
# Adds new Centers Lab object-repository rows (Notes popup / Write-off /
# Operations / Billing-edit / Comment controls) to the "Attributes & Values"
# sheet, and refreshes the Centers Lab automation credentials stored on the
# "Data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes & Values")
$data = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------
# 1. New locator rows appended after the existing Centers Lab rows
#    (sheet currently ends at row 99).  Column order below intentionally
#    mirrors how the rows were authored (xpath column first on most
#    rows, name column last) so new entries line up with the source
#    workbook.
# ---------------------------------------------------------------------

# Row 100 - CentersLab Linkto
$ws.Range("A100").Value = "Centers Lab"
$ws.Range("B100").Value = "CentersLab Linkto"
$ws.Range("C100").Value = "xpath"
$ws.Range("D100").Value = "//div[input[@name='service_note_link_to']]/input[2]"

# Row 101 - CentersLab Notes Popup
$ws.Range("A101").Value = "Centers Lab"
$ws.Range("D101").Value = "//div[@class='x-window-mc']/div/div[2]/div[6]/div[1]"
$ws.Range("B101").Value = "CentersLab Notes Popup"
$ws.Range("C101").Value = "xpath"

# Row 102 - Centers Lab Operations Option
$ws.Range("A102").Value = "Centers Lab"
$ws.Range("C102").Value = "xpath"
$ws.Range("D102").Value = "//button[text()='Operations']"

# Row 103 - Centers Lab WriteOff Option
$ws.Range("A103").Value = "Centers Lab"
$ws.Range("C103").Value = "xpath"
$ws.Range("D103").Value = "//span[text()='Write off amount']"
$ws.Range("B103").Value = "Centers Lab WriteOff Option"

$ws.Range("B102").Value = "Centers Lab Operations Option"

# Row 104 - Centers Billing Edit
$ws.Range("A104").Value = "Centers Lab"
$ws.Range("C104").Value = "xpath"
$ws.Range("D104").Value = "//span[text()='Billing edit']"
$ws.Range("B104").Value = "Centers Billing Edit"

# Row 105 - Centers Comment Add Button
$ws.Range("A105").Value = "Centers Lab"
$ws.Range("C105").Value = "xpath"
$ws.Range("D105").Value = "//button[text()='Add']"
$ws.Range("B105").Value = "Centers Comment Add Button"

# ---------------------------------------------------------------------
# 2. Refresh the Centers Lab automation credentials on the "Data" sheet.
#    Doing this between the locator rows above and the remaining rows
#    below matches the order the new unique strings were introduced in
#    the workbook.
# ---------------------------------------------------------------------
$data.Range("C11").Value = "Automation"
# Leading apostrophe keeps this a forced-text entry (preserves the cell's
# existing quote-prefix / text formatting instead of Excel reinterpreting
# the new value), matching how the password cell was re-typed in place.
$data.Range("C12").Value = "'Bright@978"

# Row 106 - Centers Lab Comment Area
$ws.Range("A106").Value = "Centers Lab"
$ws.Range("C106").Value = "xpath"
$ws.Range("D106").Value = "//textarea[@name='service_note_text']"
$ws.Range("B106").Value = "Centers Lab Comment Area"

# Row 107 - Centers Lab Save Button
$ws.Range("A107").Value = "Centers Lab"
$ws.Range("C107").Value = "xpath"
$ws.Range("D107").Value = "//button[text()='Save']"
$ws.Range("B107").Value = "Centers Lab Save Button"

# Row 108 - Centers Lab Close Button
$ws.Range("A108").Value = "Centers Lab"
$ws.Range("C108").Value = "xpath"
$ws.Range("D108").Value = "//button[text()='Close']"
$ws.Range("B108").Value = "Centers Lab Close Button"

# Rows whose wrapped text spills onto a second line get a taller row,
# matching Excel's own wrap-text autofit for the longer attribute names.
$ws.Rows.Item(102).RowHeight = 28.8
$ws.Rows.Item(103).RowHeight = 28.8
$ws.Rows.Item(105).RowHeight = 28.8

# ---------------------------------------------------------------------
# 3. Restore selection / active sheet state to match the saved workbook:
#    the "Data" sheet's own last selection moves to B21, while the
#    "Attributes & Values" sheet stays the active tab, selected at B107.
# ---------------------------------------------------------------------
$data.Range("B21").Select() | Out-Null
$ws.Activate() | Out-Null
$ws.Range("B107").Select() | Out-Null
